$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the six numeric cell values on Sheet1 (A1:B3) to the new results.
$ws.Range("A1").Value = -0.021562891158830255
$ws.Range("B1").Value = -0.021525888562684184
$ws.Range("A2").Value = 0.029042946240322089
$ws.Range("B2").Value = -0.029042946257872463
$ws.Range("A3").Value = -0.043819312237455006
$ws.Range("B3").Value = 0.043819312219303096
